$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9005614640222286
$ws.Range("J2").Value = 0.9005614640222285
$ws.Range("M2").Value = 11.651608
$ws.Range("N2").Value = 34.954824
$ws.Range("O2").Value = 0.1892813629236475
$ws.Range("P2").Value = 0.1892813629236474
$ws.Range("Q2").Value = 0.9126277320773334
$ws.Range("R2").Value = 8.213649588696001
$ws.Range("S2").Value = 0.1704595013066428
$ws.Range("T2").Value = 0.1704595013066427
$ws.Range("I3").Value = 0.9005614640222286
$ws.Range("J3").Value = 0.9005614640222285
$ws.Range("O3").Value = 0.4419371310876561
$ws.Range("P3").Value = 0.4419371310876561
$ws.Range("S3").Value = 0.3979915497780831
$ws.Range("T3").Value = 0.397991549778083
$ws.Range("I4").Value = 0.9005614640222286
$ws.Range("J4").Value = 0.9005614640222285
$ws.Range("M4").Value = 8.657178999999999
$ws.Range("N4").Value = 25.971537
$ws.Range("O4").Value = 0.1406366091439035
$ws.Range("P4").Value = 0.1406366091439035
$ws.Range("Q4").Value = 0.6780850880803333
$ws.Range("R4").Value = 6.102765792722999
$ws.Range("S4").Value = 0.1266519106257557
$ws.Range("T4").Value = 0.1266519106257556
$ws.Range("I5").Value = 0.9005614640222286
$ws.Range("J5").Value = 0.9005614640222285
$ws.Range("M5").Value = 5.488499666666667
$ws.Range("N5").Value = 16.465499
$ws.Range("O5").Value = 0.08916114387925267
$ws.Range("P5").Value = 0.08916114387925267
$ws.Range("Q5").Value = 0.4298940543912222
$ws.Range("R5").Value = 3.869046489521
$ws.Range("S5").Value = 0.08029509026579634
$ws.Range("T5").Value = 0.08029509026579634
$ws.Range("I6").Value = 0.9005614640222286
$ws.Range("J6").Value = 0.9005614640222285
$ws.Range("M6").Value = 4.091608333333333
$ws.Range("N6").Value = 12.274825
$ws.Range("O6").Value = 0.06646852536431769
$ws.Range("P6").Value = 0.06646852536431769
$ws.Range("Q6").Value = 0.3204806781861111
$ws.Range("R6").Value = 2.884326103675
$ws.Range("S6").Value = 0.05985899251348857
$ws.Range("T6").Value = 0.05985899251348856
$ws.Range("I7").Value = 0.9005614640222286
$ws.Range("J7").Value = 0.9005614640222285
$ws.Range("M7").Value = 4.463825666666667
$ws.Range("N7").Value = 13.391477
$ws.Range("O7").Value = 0.07251522760122259
$ws.Range("P7").Value = 0.07251522760122257
$ws.Range("Q7").Value = 0.3496350971092222
$ws.Range("R7").Value = 3.146715873983
$ws.Range("S7").Value = 0.06530441953246213
$ws.Range("T7").Value = 0.06530441953246212
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008648666666666667
$ws.Range("H8").Value = 0.025946
$ws.Range("I8").Value = 0.0994385359777714
$ws.Range("J8").Value = 0.09943853597777139
$ws.Range("M8").Value = 11.651608
$ws.Range("N8").Value = 34.954824
$ws.Range("O8").Value = 0.1892813629236475
$ws.Range("P8").Value = 0.1892813629236474
$ws.Range("Q8").Value = 0.1007708737226667
$ws.Range("R8").Value = 0.9069378635040001
$ws.Range("S8").Value = 0.01882186161700472
$ws.Range("T8").Value = 0.01882186161700472
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008648666666666667
$ws.Range("H9").Value = 0.025946
$ws.Range("I9").Value = 0.0994385359777714
$ws.Range("J9").Value = 0.09943853597777139
$ws.Range("O9").Value = 0.4419371310876561
$ws.Range("P9").Value = 0.4419371310876561
$ws.Range("Q9").Value = 0.2352814357542222
$ws.Range("R9").Value = 2.117532921788
$ws.Range("S9").Value = 0.04394558130957296
$ws.Range("T9").Value = 0.04394558130957296
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008648666666666667
$ws.Range("H10").Value = 0.025946
$ws.Range("I10").Value = 0.0994385359777714
$ws.Range("J10").Value = 0.09943853597777139
$ws.Range("M10").Value = 8.657178999999999
$ws.Range("N10").Value = 25.971537
$ws.Range("O10").Value = 0.1406366091439035
$ws.Range("P10").Value = 0.1406366091439035
$ws.Range("Q10").Value = 0.07487305544466667
$ws.Range("R10").Value = 0.673857499002
$ws.Range("S10").Value = 0.01398469851814782
$ws.Range("T10").Value = 0.01398469851814782
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.008648666666666667
$ws.Range("H11").Value = 0.025946
$ws.Range("I11").Value = 0.0994385359777714
$ws.Range("J11").Value = 0.09943853597777139
$ws.Range("M11").Value = 5.488499666666667
$ws.Range("N11").Value = 16.465499
$ws.Range("O11").Value = 0.08916114387925267
$ws.Range("P11").Value = 0.08916114387925267
$ws.Range("Q11").Value = 0.04746820411711111
$ws.Range("R11").Value = 0.427213837054
$ws.Range("S11").Value = 0.008866053613456319
$ws.Range("T11").Value = 0.008866053613456317
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.008648666666666667
$ws.Range("H12").Value = 0.025946
$ws.Range("I12").Value = 0.0994385359777714
$ws.Range("J12").Value = 0.09943853597777139
$ws.Range("M12").Value = 4.091608333333333
$ws.Range("N12").Value = 12.274825
$ws.Range("O12").Value = 0.06646852536431769
$ws.Range("P12").Value = 0.06646852536431769
$ws.Range("Q12").Value = 0.03538695660555556
$ws.Range("R12").Value = 0.31848260945
$ws.Range("S12").Value = 0.006609532850829115
$ws.Range("T12").Value = 0.006609532850829114
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.008648666666666667
$ws.Range("H13").Value = 0.025946
$ws.Range("I13").Value = 0.0994385359777714
$ws.Range("J13").Value = 0.09943853597777139
$ws.Range("M13").Value = 4.463825666666667
$ws.Range("N13").Value = 13.391477
$ws.Range("O13").Value = 0.07251522760122259
$ws.Range("P13").Value = 0.07251522760122257
$ws.Range("Q13").Value = 0.03860614024911112
$ws.Range("R13").Value = 0.347455262242
$ws.Range("S13").Value = 0.007210808068760454
$ws.Range("T13").Value = 0.007210808068760451
